# Actualizacion a 22 de Abril de 2020
# Adds four new daily rows (22-25 Apr 2020 / dia 47-50) to both the
# "Hoja1" (cumulative cases) and "Hoja2" (daily new cases) sheets.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Hoja1 - cumulative totals per region
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Hoja1")

$hoja1Rows = @(
    @(43939,47,145, 90,249,13,68,381,5381,55,302,678,616, 972,157,421,7,553,10088),
    @(43940,48,161, 93,263,13,69,388,5643,56,304,687,626,1045,162,424,7,566,10507),
    @(43941,49,170,104,291,13,69,403,5788,57,316,694,636,1092,167,432,7,593,10832),
    @(43942,50,180,114,326,13,70,421,6083,59,328,703,658,1113,171,445,7,605,11296)
)

$startRow = 48
for ($i = 0; $i -lt $hoja1Rows.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $hoja1Rows[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws1.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
# Column A holds dates - keep the same DD/MM/YY display format as the rest of the sheet
$ws1.Range("A48:A51").NumberFormat = "DD/MM/YY"

# -----------------------------------------------------------------
# Hoja2 - new cases per region
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Hoja2")

$hoja2Rows = @(
    @(43939,47,2,0,1,0,0,5,63,0, 8,13,3,23,3,6,0,6,133),
    @(43940,48,2,0,1,0,0,5,67,0,10,13,3,23,3,6,0,6,139),
    @(43941,49,2,0,1,0,0,6,69,1,10,13,4,25,3,6,0,7,147),
    @(43942,50,2,0,2,0,0,7,78,1,10,13,5,26,3,6,0,7,160)
)

for ($i = 0; $i -lt $hoja2Rows.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $hoja2Rows[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws2.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
$ws2.Range("A48:A51").NumberFormat = "DD/MM/YY"

# -----------------------------------------------------------------
# Update the active selection to reflect where the author ended up
# -----------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A56:AA62").Select()

$ws2.Activate()
$ws2.Range("T49").Select()

$ws1.Activate()
$ws1.Range("A56:AA62").Select()
